# "add amogus to dataBase"
# The workbook has a two-column data table on sheet "data": column A holds
# file names, column B holds the associated caption/keywords text. Row 4's
# caption cell (B4) currently holds a placeholder single space; replace it
# with the new "amogus" caption text, then update the row height (the cell
# uses a wrap-text style, so Excel grows the row to fit the new multi-line
# text) and move the active selection onto the edited cell, matching how a
# user would have made this edit interactively.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "амогус раскраска на изи амонг ас амонгас amogus among us amongus lol лол impostor импостер веселые картинки круто новинка kek кек"

# Grow row 4 to fit the new wrapped text (matches the height Excel computes
# for the other long-caption rows, e.g. row 2 which also wraps to 115.2pt).
$ws.Rows.Item(4).RowHeight = 115.2

$ws.Range("B4").Select()
